# Auto-generated Excel COM-interop script
# Applies the numeric corrections described in the commit diff to the
# Excalibur_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 284.66666
$ws.Range("I33").Value = 290.77777
$ws.Range("J33").Value = 266.33334
$ws.Range("K33").Value = 290.77777
$ws.Range("L33").Value = 266.33334
$ws.Range("M33").Value = -61.77776999999998
$ws.Range("N33").Value = -724.33334
$ws.Range("H74").Value = 8185.7617
$ws.Range("I74").Value = 4833.6665
$ws.Range("K74").Value = 4833.6665
$ws.Range("M74").Value = -3897.6665
$ws.Range("H77").Value = 8185.7617
$ws.Range("I77").Value = 4833.6665
$ws.Range("K77").Value = 24168.3325
$ws.Range("M77").Value = -19488.3325
$ws.Range("H86").Value = 2199.6
$ws.Range("I86").Value = 2199.6
$ws.Range("K86").Value = 2199.6
$ws.Range("M86").Value = -1076.6
$ws.Range("H89").Value = 2199.6
$ws.Range("I89").Value = 2199.6
$ws.Range("K89").Value = 10998
$ws.Range("M89").Value = -5382
$ws.Range("H98").Value = 971.93335
$ws.Range("I98").Value = 944.8461
$ws.Range("J98").Value = 1148
$ws.Range("K98").Value = 944.8461
$ws.Range("L98").Value = 1148
$ws.Range("M98").Value = 553.1539
$ws.Range("N98").Value = -4144
$ws.Range("H122").Value = 971.93335
$ws.Range("I122").Value = 944.8461
$ws.Range("J122").Value = 1148
$ws.Range("K122").Value = 2834.5383
$ws.Range("L122").Value = 3444
$ws.Range("M122").Value = -384.5383000000002
$ws.Range("N122").Value = -8344
$ws.Range("H137").Value = 46931096
$ws.Range("I137").Value = 83334140
$ws.Range("J137").Value = 3247433
$ws.Range("K137").Value = 250002420
$ws.Range("L137").Value = 9742299
$ws.Range("M137").Value = -249999870
$ws.Range("N137").Value = -9747399
$ws.Range("H138").Value = 11708.777
$ws.Range("J138").Value = 3189.55
$ws.Range("L138").Value = 9568.650000000001
$ws.Range("N138").Value = -19848.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7247259
$ws.Range("I32").Value = 7813264.5
$ws.Range("K32").Value = 7813264.5
$ws.Range("M32").Value = -7812977.5
$ws.Range("H61").Value = 1391937.1
$ws.Range("I61").Value = 1854190
$ws.Range("K61").Value = 1854190
$ws.Range("M61").Value = -1853978
$ws.Range("H74").Value = 2120687.8
$ws.Range("I74").Value = 2908097.2
$ws.Range("K74").Value = 2908097.2
$ws.Range("M74").Value = -2907223.2
$ws.Range("H77").Value = 2120687.8
$ws.Range("I77").Value = 2908097.2
$ws.Range("K77").Value = 14540486
$ws.Range("M77").Value = -14536118
$ws.Range("H110").Value = 2118.3
$ws.Range("I110").Value = 2236.8
$ws.Range("J110").Value = 1999.8
$ws.Range("K110").Value = 2236.8
$ws.Range("L110").Value = 1999.8
$ws.Range("M110").Value = -191.8000000000002
$ws.Range("N110").Value = -6089.8
$ws.Range("H122").Value = 4878.8
$ws.Range("I122").Value = 4872.909
$ws.Range("K122").Value = 14618.727
$ws.Range("M122").Value = -12168.727
$ws.Range("H132").Value = 449681.9
$ws.Range("I132").Value = 594451.4
$ws.Range("K132").Value = 1783354.2
$ws.Range("M132").Value = -1780824.2
$ws.Range("H136").Value = 1391937.1
$ws.Range("I136").Value = 1854190
$ws.Range("K136").Value = 5562570
$ws.Range("M136").Value = -5560020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1678
$ws.Range("I86").Value = 1531.5294
$ws.Range("K86").Value = 1531.5294
$ws.Range("M86").Value = -408.5293999999999
$ws.Range("H89").Value = 1678
$ws.Range("I89").Value = 1531.5294
$ws.Range("K89").Value = 7657.646999999999
$ws.Range("M89").Value = -2041.646999999999
$ws.Range("H99").Value = 9020.833000000001
$ws.Range("I99").Value = 3976.3845
$ws.Range("K99").Value = 3976.3845
$ws.Range("M99").Value = -2478.3845
$ws.Range("H103").Value = 59308.832
$ws.Range("J103").Value = 59308.832
$ws.Range("L103").Value = 59308.832
$ws.Range("N103").Value = -61652.832
$ws.Range("H134").Value = 482006.88
$ws.Range("I134").Value = 620825.2
$ws.Range("J134").Value = 224201.42
$ws.Range("K134").Value = 1862475.6
$ws.Range("L134").Value = 672604.26
$ws.Range("M134").Value = -1859940.6
$ws.Range("N134").Value = -677674.26

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 965.6667
$ws.Range("I25").Value = 948.5
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 948.5
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -774.5
$ws.Range("N25").Value = -1348
$ws.Range("H88").Value = 40253.168
$ws.Range("I88").Value = 36356.332
$ws.Range("J88").Value = 44150
$ws.Range("K88").Value = 36356.332
$ws.Range("L88").Value = 44150
$ws.Range("M88").Value = -35950.332
$ws.Range("N88").Value = -44962
$ws.Range("H91").Value = 40253.168
$ws.Range("I91").Value = 36356.332
$ws.Range("J91").Value = 44150
$ws.Range("K91").Value = 36356.332
$ws.Range("L91").Value = 44150
$ws.Range("M91").Value = -34952.332
$ws.Range("N91").Value = -46958
$ws.Range("H99").Value = 2037.909
$ws.Range("I99").Value = 1861.8334
$ws.Range("J99").Value = 2249.2
$ws.Range("K99").Value = 1861.8334
$ws.Range("L99").Value = 2249.2
$ws.Range("M99").Value = -363.8334
$ws.Range("N99").Value = -5245.2
$ws.Range("H122").Value = 1961.5
$ws.Range("I122").Value = 1065.8889
$ws.Range("J122").Value = 3113
$ws.Range("K122").Value = 3197.6667
$ws.Range("L122").Value = 9339
$ws.Range("M122").Value = -747.6666999999998
$ws.Range("N122").Value = -14239
$ws.Range("H126").Value = 2037.909
$ws.Range("I126").Value = 1861.8334
$ws.Range("J126").Value = 2249.2
$ws.Range("K126").Value = 5585.5002
$ws.Range("L126").Value = 6747.599999999999
$ws.Range("M126").Value = -3115.5002
$ws.Range("N126").Value = -11687.6
$ws.Range("H132").Value = 18947640
$ws.Range("I132").Value = 22737488
$ws.Range("J132").Value = 11367944
$ws.Range("K132").Value = 68212464
$ws.Range("L132").Value = 34103832
$ws.Range("M132").Value = -68209934
$ws.Range("N132").Value = -34108892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 410
$ws.Range("I8").Value = 410
$ws.Range("K8").Value = 1230
$ws.Range("M8").Value = -1091
$ws.Range("H23").Value = 377.4375
$ws.Range("J23").Value = 634.75
$ws.Range("L23").Value = 1904.25
$ws.Range("N23").Value = -2374.25
$ws.Range("H37").Value = 96999.5
$ws.Range("J37").Value = 96999.5
$ws.Range("L37").Value = 290998.5
$ws.Range("N37").Value = -291222.5
$ws.Range("H47").Value = 384.42856
$ws.Range("I47").Value = 438.2
$ws.Range("J47").Value = 250
$ws.Range("K47").Value = 1314.6
$ws.Range("L47").Value = 750
$ws.Range("M47").Value = -883.5999999999999
$ws.Range("N47").Value = -1612
$ws.Range("H132").Value = 2480.2
$ws.Range("I132").Value = 1600.25
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 14402.25
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -11872.25
$ws.Range("N132").Value = -59060
$ws.Range("H138").Value = 3666.6667
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3666.6667
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 11000.0001
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -21280.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 36461.75
$ws.Range("J39").Value = 42999
$ws.Range("L39").Value = 42999
$ws.Range("N39").Value = -44063
$ws.Range("H59").Value = 15002.632
$ws.Range("I59").Value = 15002.703
$ws.Range("J59").Value = 15000
$ws.Range("K59").Value = 15002.703
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = -14419.703
$ws.Range("N59").Value = -16166
$ws.Range("H97").Value = 1509.7037
$ws.Range("I97").Value = 1675
$ws.Range("K97").Value = 1675
$ws.Range("M97").Value = -1179
$ws.Range("H102").Value = 9997.5
$ws.Range("I102").Value = 9997.5
$ws.Range("K102").Value = 9997.5
$ws.Range("M102").Value = -8375.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 547.6
$ws.Range("I16").Value = 365.75
$ws.Range("J16").Value = 1275
$ws.Range("K16").Value = 365.75
$ws.Range("L16").Value = 1275
$ws.Range("M16").Value = -195.75
$ws.Range("N16").Value = -1615
$ws.Range("H46").Value = 1778.7
$ws.Range("I46").Value = 1026.4
$ws.Range("J46").Value = 2029.4667
$ws.Range("K46").Value = 1026.4
$ws.Range("L46").Value = 2029.4667
$ws.Range("M46").Value = -838.4000000000001
$ws.Range("N46").Value = -2405.4667
$ws.Range("H88").Value = 30000
$ws.Range("I88").Value = 30000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 30000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -29572
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 30000
$ws.Range("I91").Value = 30000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 30000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -28518
$ws.Range("N91").ClearContents()
$ws.Range("H100").Value = 9523.385
$ws.Range("I100").Value = 2042
$ws.Range("J100").Value = 34461.332
$ws.Range("K100").Value = 2042
$ws.Range("L100").Value = 34461.332
$ws.Range("M100").Value = -1501
$ws.Range("N100").Value = -35543.332
$ws.Range("H104").Value = 95260.125
$ws.Range("J104").Value = 95260.125
$ws.Range("L104").Value = 95260.125
$ws.Range("N104").Value = -102248.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2757.4878
$ws.Range("I122").Value = 2459.919
$ws.Range("K122").Value = 7379.757
$ws.Range("M122").Value = -4929.757

$wb.Save()
